# Add "PMID" column to the "studies" sheet and "notes" column to the
# "counts" sheet (new trailing header cells), then leave the selection /
# active sheet state matching the target: "counts" becomes the active
# (last-selected) sheet, with H2 selected on "studies" and F2 selected on
# "counts".

$wb = $excel.ActiveWorkbook

# --- studies sheet: add new header "PMID" in column H (row 1) ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("H1").Value = "PMID"

# --- counts sheet: add new header "notes" in column F (row 1) ---
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("F1").Value = "notes"

# --- update selections on each sheet ---
$wsStudies.Select() | Out-Null
$wsStudies.Range("H2").Select() | Out-Null

# "counts" ends up the active/selected sheet
$wsCounts.Select() | Out-Null
$wsCounts.Range("F2").Select() | Out-Null
